# BUG: Don't extract header names if none specified (gh-11733)
#
# Add a new worksheet "index_col_none" to the testmultiindex workbook. It
# holds a 2-row multiindex header ("A"/"A"/"B"/"B" over "key"/"val"/"key"/"val")
# with two data rows of 1,2,3,4 — the regression fixture for
# read_excel(..., header=[0, 1]) when no index/header names are present.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the current last sheet so it lands at the end of
# the tab strip (Worksheets.Add() with no args inserts before the active
# sheet, which is not what we want here).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "index_col_none"

# Row 1: top header level
$newSheet.Range("A1").Value = "A"
$newSheet.Range("B1").Value = "A"
$newSheet.Range("C1").Value = "B"
$newSheet.Range("D1").Value = "B"

# Row 2: second header level
$newSheet.Range("A2").Value = "key"
$newSheet.Range("B2").Value = "val"
$newSheet.Range("C2").Value = "key"
$newSheet.Range("D2").Value = "val"

# Rows 3-4: data
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = 2
$newSheet.Range("C3").Value = 3
$newSheet.Range("D3").Value = 4

$newSheet.Range("A4").Value = 1
$newSheet.Range("B4").Value = 2
$newSheet.Range("C4").Value = 3
$newSheet.Range("D4").Value = 4

# Center-align the data rows first (creates the plain "center" style), then
# bold + center-align the header rows (reuses the bold/center header style
# already used elsewhere in this workbook) so the style table comes out
# matching the original authoring order.
$dataRange = $newSheet.Range("A3:D4")
$dataRange.HorizontalAlignment = -4108  # xlCenter

$headerRange = $newSheet.Range("A1:D2")
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.Font.Bold = $true

# Match the saved cursor position/tab seen in the authored file.
$newSheet.Range("G23").Select()
